# [MOSIP-14369] Fix: boolean values
#
# The J column ("is_active") cells J2:J73 were storing the boolean result of
# the formula =TRUE() (numeric 1 under the hood). They must instead hold the
# literal text string "TRUE". A plain Value assignment of "TRUE" gets
# auto-coerced back to a boolean by the engine's type inference, so we stage
# the literal text in a scratch cell (using a leading apostrophe to force
# text interpretation) and then copy/paste-special (values-only) it onto
# every target cell - this preserves each cell's existing style/number
# format while only replacing the stored value/type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z1")
$scratch.Value = "'TRUE"
$scratch.Copy()

for ($r = 2; $r -le 73; $r++) {
    $ws.Cells.Item($r, 10).PasteSpecial(-4163, 0, $false, $false)  # xlPasteValues
}

$scratch.Clear()

# Restore the selection/scroll state to match the committed view: the sheet
# is scrolled so row 51 is the top row, and J2:J73 is selected with J2 active.
$ws.Range("J2:J73").Select()
$win = $excel.ActiveWindow()
$win.ScrollRow = 51
$win.ScrollColumn = 1
